$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new data row at row 64, pushing existing rows 64-98 down to 65-99.
$ws.Rows.Item(64).Insert()

# Populate the new row 64 with the new market-record values.
$ws.Range("A64").Value2 = 2
$ws.Range("B64").Value2 = "Comercializadora del Agro de Limarí"
$ws.Range("C64").Value2 = "Coquimbo"
$ws.Range("D64").Value2 = 44553
$ws.Range("E64").Value2 = 4
$ws.Range("F64").Value2 = "Fruta"
$ws.Range("G64").Value2 = 100109
$ws.Range("H64").Value2 = "Uva"
$ws.Range("I64").Value2 = 100109001
$ws.Range("J64").Value2 = "Uva"
$ws.Range("K64").Value2 = "Flame Seedless"
$ws.Range("L64").Value2 = "Primera"
$ws.Range("M64").Value2 = 600
$ws.Range("N64").Value2 = 5500
$ws.Range("O64").Value2 = 6000
$ws.Range("P64").Value2 = 5750
$ws.Range("Q64").Value2 = "$/bandeja 10 kilos"
$ws.Range("R64").Value2 = "Provincia de Limarí"
$ws.Range("S64").Value2 = 575
$ws.Range("T64").Value2 = 10
